$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the existing
# header cells (e.g. G1: bold font, thin border, centered/top alignment) by
# copying G1's format onto H1 instead of creating a brand-new style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Header text for the new column.
$ws.Range("H1").Value = "Save"

# New data value for the single data row.
$ws.Range("H2").Value = 0
